$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded. Insert a new row at 149 so the
# existing rows 149:163 (and their formatting) shift down to 150:164, then
# populate the new row with the latest week's data (same market / category
# metadata as the surrounding rows).
$ws.Rows.Item(149).Insert()

$ws.Range("A149").Value = 5
$ws.Range("B149").Value = "Macroferia Regional de Talca"
$ws.Range("C149").Value = "Maule"
$ws.Range("D149").Value = 44449
$ws.Range("E149").Value = 7
$ws.Range("F149").Value = 100112009
$ws.Range("G149").Value = "Acelga"
$ws.Range("H149").Value = "Sin especificar"
$ws.Range("I149").Value = "Primera"
$ws.Range("J149").Value = 500
$ws.Range("K149").Value = 2500
$ws.Range("L149").Value = 2500
$ws.Range("M149").Value = 2500
$ws.Range("N149").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O149").Value = "Región del Maule"
$ws.Range("P149").Value = 625
$ws.Range("Q149").Value = 4
$ws.Range("R149").Value = "Hortaliza"
